# TICKET4769: template update to satisfy latest request
# Insert a new "Number of Rooms" row right after the existing "Room Nights"
# row (row 20), pushing "Flex Minimum/Max" and everything below it down by
# one row. Excel copies the formatting of the row above on insert, so this
# naturally reproduces the same styles/merged-cell layout as the rows it
# shifted from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20:24 down to 21:25, inserting a blank row 20 formatted like
# row 19 above it.
$ws.Rows.Item(20).Insert()

# Populate the new row's label/placeholder cells.
$ws.Cells.Item(20, 1).Value = "Number of Rooms"
$ws.Cells.Item(20, 2).Value = "%%PACKAGE_NUM_ROOMS%%"

# Make sure the new row keeps the same explicit row height as its
# neighbours (15pt, custom height) instead of the default autosize height.
$ws.Rows.Item(20).RowHeight = 15
